$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20-24 down to 21-25
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with data
$ws.Cells.Item(20, 1).Value = 4
$ws.Cells.Item(20, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(20, 3).Value = "Los Lagos"
$ws.Cells.Item(20, 4).Value = 44736
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value = 10
$ws.Cells.Item(20, 6).Value = 100112012
$ws.Cells.Item(20, 7).Value = "Espinaca"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 35
$ws.Cells.Item(20, 11).Value = 13000
$ws.Cells.Item(20, 12).Value = 13000
$ws.Cells.Item(20, 13).Value = 13000
$ws.Cells.Item(20, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 1300
$ws.Cells.Item(20, 17).Value = 10
$ws.Cells.Item(20, 18).Value = "Hortaliza"
